# BIS-1002: remove "Internal Assignment" column (column O) content from
# the sample-type export/import test fixture. Clearing these cells drops
# the last reference to the "Internal Assignment" shared string, so it is
# automatically pruned from sharedStrings.xml on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4:O7").ClearContents()
